# Refactor and add support for arbitrary CSV delimiters.
# Update the example fixture workbook (Sheet2 header row) so the sample
# data exercises a custom delimiter / quoting edge cases: a value that
# starts with a digit, an all-caps value, and a value containing an
# underscore, a question mark and a typographic apostrophe.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A1").Value = "2dd"
$ws2.Range("B1").Value = "EEE"
$ws2.Range("C1").Value = "fff ggg_?h" + [char]0x2019 + "i"

# Leave the selection on B1, matching the saved view state.
$ws2.Activate() | Out-Null
$ws2.Range("B1").Select() | Out-Null
